$d = $word.ActiveDocument

# -----------------------------------------------------------------------
# Change 1: fix typo "chứ" -> "chứa" inside the RESTful/EXP paragraph.
# -----------------------------------------------------------------------
$found1 = $d.Content.Find.Execute(
    "cú pháp JSON chứ những thông số", $true, $false, $false, $false, $false,
    $true, 1, $false, "cú pháp JSON chứa những thông số", 2)

# -----------------------------------------------------------------------
# Change 2: "Tự động cấu hình mỗi subnet (" + "PAC)" -> merge into a single
# run "Tự động cấu hình mỗi subnet (PAC)", and relocate the "_GoBack"
# bookmark (that used to sit between those two runs) to the very end of
# the following paragraph ("Như đã thấy ở hình 1...").
# -----------------------------------------------------------------------

# The "_GoBack" bookmark currently sits right between the two runs that
# are about to be merged; remove it from there first.
$goBack = $d.Bookmarks("_GoBack")
$goBack.Delete()

# Re-running Find/Replace over the full heading text (with no formatting
# changes requested) merges the two adjoining runs that contained
# "Tự động cấu hình mỗi subnet (" and "PAC)" into a single run, taking on
# the formatting of the first (bold+italic) run, exactly like the diff.
$found2 = $d.Content.Find.Execute(
    "Tự động cấu hình mỗi subnet (PAC)", $true, $false, $false, $false,
    $false, $true, 1, $false, "Tự động cấu hình mỗi subnet (PAC)", 2)

# Now re-insert the "_GoBack" bookmark at the end of the next paragraph,
# right after its last visible character and before the paragraph mark.
# A collapsed (zero-length) Range placed directly at that boundary is
# mishandled by this runtime, so we work around it: insert a one-character
# placeholder there, wrap a bookmark around that character, then delete
# the placeholder again - the bookmark naturally collapses to the
# now-empty spot, which is exactly the position we need.
$p34 = $d.Paragraphs(34)
$insertPos = $p34.Range.End - 1
$placeholder = $d.Range($insertPos, $insertPos)
$placeholder.InsertAfter("X")

$bmRange = $d.Range($insertPos, $insertPos + 1)
$d.Bookmarks.Add("_GoBack", $bmRange)

$cleanup = $d.Range($insertPos, $insertPos + 1)
$cleanup.Text = ""
